$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarterly data (D:K) to (F:M)
$ws.Columns("D:E").Insert()

# Copy number formats/styles from the (now-shifted) old column D -- now column F -- onto the two new columns.
# Done per contiguous data block so the untouched label-only rows (5, 6, 37, 79) are not given
# spurious empty D/E cells.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarters of data (column D = Q4 2018, column E = Q3 2018)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 69700
$ws.Range("E8").Value = 64800
$ws.Range("D9").Value = 23000
$ws.Range("E9").Value = 18900
$ws.Range("D10").Value = 46700
$ws.Range("E10").Value = 45900
$ws.Range("D12").Value = 15000
$ws.Range("E12").Value = 13100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 78300
$ws.Range("E17").Value = 70000
$ws.Range("D18").Value = -8600
$ws.Range("E18").Value = -5200
$ws.Range("D20").Value = -2200
$ws.Range("E20").Value = -2400
$ws.Range("D21").Value = -8600
$ws.Range("E21").Value = -5300
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = -10800
$ws.Range("E23").Value = -7600
$ws.Range("D24").Value = 4800
$ws.Range("E24").Value = 4100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -15600
$ws.Range("E26").Value = -11700
$ws.Range("D27").Value = -15600
$ws.Range("E27").Value = -11700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2200
$ws.Range("E32").Value = 2400
$ws.Range("D33").Value = -15600
$ws.Range("E33").Value = -11700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -15600
$ws.Range("E35").Value = -11700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 140600
$ws.Range("E41").Value = 123000
$ws.Range("D42").Value = 105300
$ws.Range("E42").Value = 104700
$ws.Range("D43").Value = 50800
$ws.Range("E43").Value = 45300
$ws.Range("D44").Value = 22600
$ws.Range("E44").Value = 21600
$ws.Range("D45").Value = 2100
$ws.Range("E45").Value = 2200
$ws.Range("D46").Value = 321400
$ws.Range("E46").Value = 296800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 15400
$ws.Range("E48").Value = 15900
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3100
$ws.Range("E52").Value = 2800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 339800
$ws.Range("E54").Value = 315500
$ws.Range("D57").Value = 64600
$ws.Range("E57").Value = 48100
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 64600
$ws.Range("E60").Value = 48100
$ws.Range("D61").Value = 149300
$ws.Range("E61").Value = 149200
$ws.Range("D62").Value = 13700
$ws.Range("E62").Value = 3300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 227500
$ws.Range("E66").Value = 200600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -643700
$ws.Range("E72").Value = -628000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 112300
$ws.Range("E76").Value = 114900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -15600
$ws.Range("E81").Value = -11700
$ws.Range("D83").Value = 2200
$ws.Range("E83").Value = 2300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 16800
$ws.Range("E89").Value = 5600
$ws.Range("D91").Value = -800
$ws.Range("E91").Value = -600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -1800
$ws.Range("E94").Value = -1000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 2600
$ws.Range("E100").Value = 3900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 17600
$ws.Range("E102").Value = 8500
